$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) stays the same content; nothing to change there ---

# --- Wipe any stray formatting left over in the old A2:E5 block (e.g. the
#     leftover underline style parked on B5) before repopulating it ---
$ws.Range("A2:E7").ClearFormats()

# --- Make sure column E (FECHANAC) keeps its "text" number format for all
#     data rows, so the dd-mm-yyyy-looking strings are not auto-converted
#     into date serials. ---
$ws.Range("E2:E7").NumberFormat = "@"

# --- Row 2 : AUT_ODON_001 ---
$ws.Range("A2").Value = "AUT_ODON_001"
$ws.Range("B2").Value = "Nancy"
$ws.Range("C2").Value = "Meza"
$ws.Range("D2").Value = "Zuñiga"
$ws.Range("E2").Value = "01-02-1985"

# --- Row 3 : AUT_ODON_002 ---
$ws.Range("A3").Value = "AUT_ODON_002"
$ws.Range("B3").Value = "Jose"
$ws.Range("C3").Value = "Pérez"
$ws.Range("D3").Value = "Gónzales"
$ws.Range("E3").Value = "15-12-1993"

# --- Row 4 : AUT_PSIC_001 (new row) ---
$ws.Range("A4").Value = "AUT_PSIC_001"
$ws.Range("B4").Value = "Marlenne"
$ws.Range("C4").Value = "Cruz"
$ws.Range("D4").Value = "Perez"
$ws.Range("E4").Value = "01-02-1985"

# --- Row 5 : AUT_PSIC_002 (new row) ---
$ws.Range("A5").Value = "AUT_PSIC_002"
$ws.Range("B5").Value = "Mario"
$ws.Range("C5").Value = "Sebastian"
$ws.Range("D5").Value = "Cerro"
$ws.Range("E5").Value = "15-12-1993"

# --- Row 6 : AUT_NUTRI_001 (new row) ---
$ws.Range("A6").Value = "AUT_NUTRI_001"
$ws.Range("B6").Value = "Gibran"
$ws.Range("C6").Value = "Pedraza"
$ws.Range("D6").Value = "Morroy"
$ws.Range("E6").Value = "01-02-1985"

# --- Row 7 : AUT_NUTRI_002 (new row) ---
$ws.Range("A7").Value = "AUT_NUTRI_002"
$ws.Range("B7").Value = "Kaarina"
$ws.Range("C7").Value = "Jiménez"
$ws.Range("D7").Value = "López"
$ws.Range("E7").Value = "15-12-1993"

# --- Column A got wider to fit the longer "AUT_xxx_xxx" codes ---
$ws.Columns.Item(1).ColumnWidth = 21

# --- Reviewer left the selection parked outside of the table ---
$ws.Range("H5").Select()
